# Add the new "Digger" class to the Lvl 99 stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the first unused "New Class" column header (U7, merged U7:V7) to "Digger".
$ws.Range("U7").Value = "Digger"

# Enter the Lvl 99 percentage modifiers for the Digger class.
$ws.Range("U9").Value = 1.02
$ws.Range("U10").Value = 0
$ws.Range("U11").Value = 1.05
$ws.Range("U12").Value = 1.13
$ws.Range("U13").Value = 0
$ws.Range("U14").Value = 0.94
$ws.Range("U15").Value = 0.85
$ws.Range("U16").Value = 1.01

# Leave the cursor where the user ended up after the edit.
$ws.Application.Goto($ws.Range("U17"))
